$wb = $excel.ActiveWorkbook

# Overview sheet: Latest HO Xliff Generate Date
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-04 19:09:30"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-04 19:09:25"
$wsZhCn.Range("K2").Value = "2016-09-04 19:10:10"

# de-de sheet: Correspond Handoff Datetime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-04 19:10:19"
